# PvsI model-fitting update.
# The respirometry pipeline re-fit the chamber "volume" (column T) per
# sample from the calibration model, which cascades into the dependent
# rate columns (Z = rate.abs, AB = rate.a.spec, AD = rate.output).
# rate.abs (Z)      = (rate - adjustment) * volume          [N,O -> T]
# rate.a.spec (AB)  = rate.abs / area                        [Z -> V]
# rate.output (AD)  = rate.a.spec                             [AB]

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fitted chamber volumes (L), keyed by data row. Rows 2-9 are the
# Respiration (dark) block, rows 10-17 the Photosynthesis (light) block;
# each sample (Gen01..Gen08) keeps the same fitted volume across both.
$newVolume = @{
    2  = 0.1405756097560976
    3  = 0.144009756097561
    4  = 0.1462634146341464
    5  = 0.1500390243902439
    6  = 0.1463707317073171
    7  = 0.1465658536585366
    8  = 0.1449658536585366
    9  = 0.1544
    10 = 0.1405756097560976
    11 = 0.144009756097561
    12 = 0.1462634146341464
    13 = 0.1500390243902439
    14 = 0.1463707317073171
    15 = 0.1465658536585366
    16 = 0.1449658536585366
    17 = 0.1544
}

foreach ($row in ($newVolume.Keys | Sort-Object)) {
    $volume = $newVolume[$row]
    $ws.Range("T$row").Value = $volume

    $rate = $ws.Range("N$row").Value()
    $adjustment = $ws.Range("O$row").Value()
    $area = $ws.Range("V$row").Value()

    $rateAbs = ($rate - $adjustment) * $volume
    $ws.Range("Z$row").Value = $rateAbs

    # Rows with zero chamber area have no meaningful area-specific rate
    # (matches the pre-existing blank/"Inf" handling in those rows).
    if ($area -ne 0) {
        $rateASpec = $rateAbs / $area
        $ws.Range("AB$row").Value = $rateASpec
        $ws.Range("AD$row").Value = $rateASpec
    }
}
